$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as scraped by GitHub Actions.
# Column D (Price) values must remain plain text (they may look numeric),
# so we force a text number format before assigning, then restore the
# default "Normal" style so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.148.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.724.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.531"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.129.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0749"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "190.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("B23").Value = "Monero"
$ws.Range("C23").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "143.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0480"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.878"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.127.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.520"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0155"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.791"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0115"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0522"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.417"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0928"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.86%  "
$ws.Range("E51").Value = "  -0.59%  "
